# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column A (cursos-jefe-explotacion) switches from an iaest-dimension to an
# iaest-measure; column D (municipio-nombre) switches the other way, from an
# iaest-measure to an sdmx-dimension:refArea. The now-unused mapping file
# reference for column A (row 5) is cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: cursos-jefe-explotacion -> now a measure, not a dimension
$ws.Range("A2").Value = "iaest-measure:cursos-jefe-explotacion"
$ws.Range("A3").Value = "medida"
$ws.Range("A4").Value = "xsd:int"
$ws.Range("A5").Clear()

# Column D: municipio-nombre -> now a dimension (refArea), not a measure
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"
